# Regional Availability Factor workbook update
# - About: bump the "last updated" date
# - RAF-capacity: raise the RAF for hydrogen combustion turbine / combined cycle
#   from 0.3 to 1, widen column A a bit, and leave the view parked there
#   (matches the author having been last working on that tab).

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

$wsCap = $wb.Worksheets.Item("RAF-capacity")
$wsCap.Range("B24").Value = 1
$wsCap.Range("B25").Value = 1
$wsCap.Columns.Item(1).ColumnWidth = 28.14

$wsCap.Select() | Out-Null
$wsCap.Range("A14").Select() | Out-Null
$wsCap.Range("B25").Select() | Out-Null

$win = $excel.ActiveWindow
$win.Zoom = 80
